$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update T_Boil (B5) with the new value; Boils_per_Charge (B6) formula
# will recalculate automatically from B3/(B4*B5).
$ws.Range("B5").Value = 141.79316571261199

$excel.Calculate()
